$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the existing row 6 data down into the new row 7
# (row 7 ends up with the values that used to be in row 6)
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(7, $col).Value = $ws.Cells.Item(6, $col).Value()
}
# Copy the date number format from D6 to D7 so the date format carries over
$ws.Range("D7").NumberFormat = $ws.Range("D6").NumberFormat

# Now update row 3
$ws.Cells.Item(3, 4).Value = 45251   # D3
$ws.Cells.Item(3, 13).Value = 15     # M3

# Update row 4
$ws.Cells.Item(4, 4).Value = 45239   # D4
$ws.Cells.Item(4, 13).Value = 25     # M4
$ws.Cells.Item(4, 18).Value = "Provincia de San Felipe de Aconcagua"  # R4

# Update row 5
$ws.Cells.Item(5, 4).Value = 45244   # D5
$ws.Cells.Item(5, 13).Value = 70     # M5

# Update row 6 (becomes the former "Región Metropolitana" row)
$ws.Cells.Item(6, 4).Value = 45243   # D6
$ws.Cells.Item(6, 13).Value = 50     # M6
$ws.Cells.Item(6, 14).Value = 35000  # N6
$ws.Cells.Item(6, 15).Value = 35000  # O6
$ws.Cells.Item(6, 16).Value = 35000  # P6
$ws.Cells.Item(6, 18).Value = "Región Metropolitana"  # R6
$ws.Cells.Item(6, 19).Value = 7000   # S6
